$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text values look like plain numbers need an explicit
# Text number format first, otherwise Excel auto-converts the assigned
# string into a numeric value (losing e.g. trailing zeros).
$ws.Range("D2").Value = "40.832.22"
$ws.Range("E2").Value = "  -4.47%  "
$ws.Range("D3").Value = "2.446.99"
$ws.Range("E3").Value = "  -3.75%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.38"
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.82"
$ws.Range("E6").Value = "  -7.53%  "
$ws.Range("E7").Value = "  -4.22%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.499"
$ws.Range("E9").Value = "  -5.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.36"
$ws.Range("E10").Value = "  -8.34%  "
$ws.Range("E11").Value = "  -3.46%  "
$ws.Range("E12").Value = "  -0.76%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.91"
$ws.Range("E13").Value = "  -6.19%  "
$ws.Range("D14").Value = "2.823.51"
$ws.Range("E14").Value = "  -3.59%  "
$ws.Range("D15").Value = "2.449.14"
$ws.Range("E15").Value = "  -4.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.36"
$ws.Range("E16").Value = "  -8.94%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.782"
$ws.Range("E17").Value = "  -3.87%  "
$ws.Range("D18").Value = "40.872.06"
$ws.Range("E18").Value = "  -4.34%  "
$ws.Range("D20").Value = "0.0₃0910"
$ws.Range("E20").Value = "  -4.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.49"
$ws.Range("E21").Value = "  -6.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "66.79"
$ws.Range("E22").Value = "  -4.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.70"
$ws.Range("E24").Value = "  -4.97%  "
$ws.Range("E25").Value = "  -7.06%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.42"
$ws.Range("E27").Value = "  -6.75%  "
$ws.Range("E28").Value = "  -4.32%  "
$ws.Range("E29").Value = "  -5.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.01"
$ws.Range("E30").Value = "  -8.75%  "
$ws.Range("E31").Value = "  -2.46%  "
$ws.Range("E32").Value = "  -4.46%  "
$ws.Range("E33").Value = "  -1.07%  "
$ws.Range("E34").Value = "  -9.01%  "
$ws.Range("E35").Value = "  -6.05%  "
$ws.Range("E36").Value = "  -5.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.06"
$ws.Range("E37").Value = "  -6.76%  "
$ws.Range("E38").Value = "  -7.81%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.113"
$ws.Range("E39").Value = "  -4.67%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.102"
$ws.Range("E40").Value = "  -8.98%  "
$ws.Range("E41").Value = "  -4.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.08"
$ws.Range("E42").Value = "  -5.01%  "
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("D44").Value = "1.959.21"
$ws.Range("E44").Value = "  -1.27%  "
$ws.Range("E45").Value = "  -5.48%  "
$ws.Range("E46").Value = "  -9.08%  "
$ws.Range("E47").Value = "  -2.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "76.06"
$ws.Range("E48").Value = "  -6.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "69.17"
$ws.Range("E49").Value = "  -5.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "97.01"
$ws.Range("E50").Value = "  -3.90%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.178"
$ws.Range("E51").Value = "  -7.25%  "
